$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the "centered" style cleanly on a single cell first (so the engine
# settles on one final style definition instead of leaving a half-applied
# transitional style behind), then copy/paste that format onto every other
# cell that needs it - this reuses the same style id everywhere.
$a1 = $ws.Cells.Item(1, 1)
$a1.HorizontalAlignment = -4108   # xlCenter
$a1.VerticalAlignment = -4108     # xlVAlignCenter

$a1.Copy()
$ws.Range("A1:C4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E7:F8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Values ---
$ws.Range("A1").Value = "First entry"
$ws.Range("D5").Value = "Second entry"
$ws.Range("E7").Value = "second merge"

# --- Merges ---
$ws.Range("A1:C4").Merge()
$ws.Range("E7:F8").Merge()
